$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Insert a new column before column H (old H:K -> new I:L) ---
# This gives room for the new "Range" column between DMG (G) and Special ability (I).
$ws.Columns.Item(8).Insert() | Out-Null
$ws.Columns.Item(8).ColumnWidth = $ws.Columns.Item(7).ColumnWidth

# --- Row 14 (Area Patrol Bot): now 2 HP, "Box, no" movement, FOV range 3, updated note ---
$ws.Range("B14").Value = 2
$ws.Range("C14").Value = "Box, no"

# --- Row 15 (Patrol Bot): now 2 HP, "Linear, no" movement, FOV range 3 ---
$ws.Range("B15").Value = 2
$ws.Range("C15").Value = "Linear, no"

$ws.Range("L14").Value = "First level, non pursue, single hit kills"

# --- Row 17 (Swarmer): now 4 HP, expanded note ---
$ws.Range("B17").Value = 4
$ws.Range("L17").Value = "Introduction to sound. Swarmers are a bit easy if they are single hit kills - just find a corridor"

# --- Header row: new "Range" column header ---
$ws.Range("H1").Value = "Range"

# --- Weapon table header: add Range / Throw range columns ---
$ws.Range("E34").Value = "Range"
$ws.Range("F34").Value = "Throw range"

# --- Row 18 held a stray "`" in A18 - clear it (row disappears once empty) ---
$ws.Range("A18").ClearContents()

# --- Row 19 (Rotating turret): FOV range 8 ---
$ws.Range("H19").Value = 8

$ws.Range("H14").Value = 3
$ws.Range("H15").Value = 3

# --- Vibroblade throw range ---
$ws.Range("F37").Value = 5

# --- Pistol range / throw range ---
$ws.Range("E38").Value = 8
$ws.Range("F38").Value = 3

# --- Shotgun range / throw range ---
$ws.Range("E39").Value = "to 10"
$ws.Range("F39").Value = 3

# --- Final selection as left by the author ---
$ws.Range("F37").Select() | Out-Null
